# "Update with a new song" - add "Little Black Submarines" below the
# existing "Uprising" entry in the band's repertoire list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Little Black Submarines"

# Resize column A so the new, longer title is fully visible (mirrors the
# author double-clicking the column border to auto-fit after typing).
$ws.Columns("A:A").AutoFit() | Out-Null

# Leave the selection where the author ended up after entering the song.
$ws.Range("A6").Select() | Out-Null
